# "delete cascade implemented and view bugs fixed"
# Mark a batch of previously-open "Problem" rows as Fixed (column C),
# and fix up the saved view/selection state on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows whose "Urgent"/status column (C) should now read "Fixed".
$fixedRows = @(8, 12, 30, 31, 32, 33, 34, 35, 37, 38)

foreach ($r in $fixedRows) {
    $ws.Cells.Item($r, 3).Value = "Fixed"
}

# View bug fixes: scroll position and active selection had drifted;
# pin the top-left visible cell and the active selection back to where
# they belong now that the list has grown.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 21
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C35").Select()
